# repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F10").Value = -8
$ws.Range("F11").Value = 13
$ws.Range("F14").Value = -3
$ws.Range("F20").Value = -3
$ws.Range("F21").Value = 5
$ws.Range("F22").Value = -2
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = -1
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = -1
$ws.Range("F34").Value = -2
$ws.Range("F35").Value = -6
